$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank/index column (column A), shifting B->A and C->B
$ws.Columns.Item(1).Delete()
